$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Sema4d", "Plxnb1", "ECs", 3, 1, 76.08252666666667, 228.24758, 0.95878149807566, 0.95878149807566, 2, 0.6666666666666666, 0.6893159999999999, 2.067948, 0.1083604551316437, 0.1083604551316437, 52.44490295175999, 472.00412656584, 0.1038939995032777, 0.1038939995032777),
    @("ECs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 76.08252666666667, 228.24758, 0.95878149807566, 0.95878149807566, 2, 0.6666666666666666, 0.082959, 0.248877, 0.01304115238477858, 0.01304115238477858, 6.31173032974, 56.80557296766001, 0.01250361562011097, 0.01250361562011097),
    @("ECs", "Sema4d", "Plxnb1", "sCs", 3, 1, 76.08252666666667, 228.24758, 0.95878149807566, 0.95878149807566, 3, 1, 5.589049333333333, 16.767148, 0.8785983924835777, 0.8785983924835776, 425.2289949446488, 3827.06095450184, 0.8423838829522714, 0.8423838829522713),
    @("FAPs", "Sema4d", "Plxnb1", "ECs", 3, 1, 0.9347409999999999, 2.804223, 0.01177947704364805, 0.01177947704364805, 2, 0.6666666666666666, 0.6893159999999999, 2.067948, 0.1083604551316437, 0.1083604551316437, 0.6443319271559999, 5.798987344404, 0.001276429493662452, 0.001276429493662451),
    @("FAPs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 0.9347409999999999, 2.804223, 0.01177947704364805, 0.01177947704364805, 2, 0.6666666666666666, 0.082959, 0.248877, 0.01304115238477858, 0.01304115238477858, 0.077545178619, 0.697906607571, 0.0001536179551392153, 0.0001536179551392153),
    @("FAPs", "Sema4d", "Plxnb1", "sCs", 3, 1, 0.9347409999999999, 2.804223, 0.01177947704364805, 0.01177947704364805, 3, 1, 5.589049333333333, 16.767148, 0.8785983924835777, 0.8785983924835776, 5.224313562889332, 47.01882206600399, 0.01034942959484638, 0.01034942959484638),
    @("sCs", "Sema4d", "Plxnb1", "ECs", 3, 1, 2.336085333333334, 7.008256, 0.02943902488069198, 0.02943902488069198, 2, 0.6666666666666666, 0.6893159999999999, 2.067948, 0.1083604551316437, 0.1083604551316437, 1.610300997632, 14.492708978688, 0.003190026134703566, 0.003190026134703566),
    @("sCs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 2.336085333333334, 7.008256, 0.02943902488069198, 0.02943902488069198, 2, 0.6666666666666666, 0.082959, 0.248877, 0.01304115238477858, 0.01304115238477858, 0.193799303168, 1.744193728512, 0.0003839188095283921, 0.0003839188095283921),
    @("sCs", "Sema4d", "Plxnb1", "sCs", 3, 1, 2.336085333333334, 7.008256, 0.02943902488069198, 0.02943902488069198, 3, 1, 5.589049333333333, 16.767148, 0.8785983924835777, 0.8785983924835776, 13.05649617487644, 117.508465573888, 0.02586507993646002, 0.02586507993646002)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
